# Wizyty.xlsx: append 3 new "visit" rows (rows 3-5) below the existing
# header (row 1) and first data row (row 2), mirroring the columns:
# A=Imie, B=Nazwisko, C=Telefon, D=Data, E=Godzina, F=Email, G=Id

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 — duplicate of row 2's patient details, new visit id "3"
$ws.Range("A3").Value = "Wiktor"
$ws.Range("B3").Value = "Kowalski"
$ws.Range("C3").Value = "'23412343"
$ws.Range("D3").Value = "19.01.2023"
$ws.Range("E3").Value = "'19.00"
$ws.Range("F3").Value = "wiktor.k.2002@icloud.com"
$ws.Range("G3").Value = "'3"

# Row 4 — new entry
$ws.Range("A4").Value = "'12"
$ws.Range("B4").Value = "'2"
$ws.Range("C4").Value = "'21323123121"
$ws.Range("D4").Value = "20.10.2023"
$ws.Range("E4").Value = "'14.11"
$ws.Range("F4").Value = "'3"
$ws.Range("G4").Value = "'4"

# Row 5 — new entry
$ws.Range("A5").Value = "XD"
$ws.Range("B5").Value = "xd"
$ws.Range("C5").Value = "'23121243312"
$ws.Range("D5").Value = "20.01.2023"
$ws.Range("E5").Value = "'15.11"
$ws.Range("F5").Value = "ddd"
$ws.Range("G5").Value = "'5"

# The values above that look numeric (ids, phone numbers, times) must stay
# plain text (shared strings), matching the source workbook where every
# cell -- regardless of content -- was written as text. Typing them via
# COM with a leading apostrophe forces text storage, but Excel also stamps
# the cell with a "quote prefix" style. Re-apply the sheet's untouched
# default style (taken from A1, which was never modified) on top so the
# cells end up as plain shared-string text with no leftover formatting,
# exactly like the rest of the sheet.
$baseStyle = $ws.Range("A1").Style
$textifiedCells = @("C3","E3","G3","A4","B4","C4","E4","F4","G4","C5","E5","G5")
foreach ($addr in $textifiedCells) {
    $ws.Range($addr).Style = $baseStyle
}
